# [PHOENIX-5860] Updated Create and search Trade License Screen
#
# - tradeOwnerDetails!B2 (aadhaarNumber): "11111111" -> "123456789123",
#   right-aligned to better fit the longer numeric-looking value.
# - tradeDetails sheet: normalize the header/data row formatting (drop the
#   stray per-row font override) and leave the F/H columns on the sheet's
#   normal/default style.
# - Selections/active tab restored to the owner-details screen (the first
#   tab of the "Create Trade License" workbook) with the cursor sitting on
#   the aadhaar column; the trade-details tab keeps its own remembered
#   selection but is no longer the active tab.

$wb = $excel.ActiveWorkbook

# --- tradeOwnerDetails -------------------------------------------------
$ws1 = $wb.Worksheets.Item("tradeOwnerDetails")

# Update the sample aadhaar number used by the functional test fixture.
$ws1.Range("B2").Value = "123456789123"

# The longer value reads better right-aligned (matches the other numeric
# columns on this sheet).
$ws1.Range("B2").HorizontalAlignment = -4152  # xlRight

# --- tradeDetails --------------------------------------------------------
$ws3 = $wb.Worksheets.Item("tradeDetails")

# Re-apply the explicit font to the header + data rows so they pick up the
# sheet's normal cell style instead of the stray duplicate style that only
# the F/H (tradeAreaWeightOfPremises / tradeCommencementDate) columns used.
$ws3.Range("A1:H1").Font.Name = "Arial"
$ws3.Range("A2:G2").Font.Name = "Arial"

# --- Selections / active tab --------------------------------------------
# Remember a selection on the trade-details tab (no longer the active tab).
$ws3.Range("F3").Select()

# Make the owner-details tab the active one again, with the cursor on the
# aadhaar number cell that was just edited.
$ws1.Activate()
$ws1.Range("C9").Select()
